# Update "Sprint 1 Example" burndown sheet for the new Day/Week 5 data
# point:
#  - Total points (C5) and total hours (C6) grow to account for the
#    added work.
#  - Day 5 (column G) effort estimate, hours-remaining and points-burned
#    inputs are filled in (previously blank / downstream "#N/A"
#    placeholders).
#  - The "Remaining" backlog series (row 24, days 1-4) is revised.
# All the other cells on the sheet (O7, rows 10-12/14/16/18/23/25, the
# chart caches, etc.) are plain formulas/derived caches that recompute
# from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1 Example")

$ws.Range("C5").Value = 433
$ws.Range("C6").Value = 52

$ws.Range("G9").Value = 5
$ws.Range("G13").Value = 122
$ws.Range("G17").Value = 52

$ws.Range("C24").Value = 41
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = 25
$ws.Range("F24").Value = 19

$wb.Application.CalculateFull()

# Work around a stale-shared-formula recalculation quirk: G10 and G12
# depend on G14 (via LOOKUP($C14:$L14) / a direct reference) which just
# flipped from #N/A to a real number. Re-asserting their formulas forces
# a fresh evaluation so they pick up the new G14 value instead of a
# cached pre-edit result.
$ws.Range("G10").Formula = "=IF(LEN(G9)>0,(LOOKUP(9.99999999999999E+307,`$C14:`$L14))/(COUNT(`$C`$14:`$L`$14))*G9,NA())"
$ws.Range("G12").Formula = "=IF(G9>0,IF(G13=0,F12,`$C5-G14),NA())"

$wb.Application.CalculateFull()

# Match the author's new active selection on this sheet.
$ws.Activate()
$ws.Range("C25").Select()
